$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Name -eq $name) {
            return $shp
        }
    }
    return $null
}

# "A" + " " + "slide" runs were merged into a single "A slide" run.
$titleShape = Get-ShapeByName $s "Title 1"
$titleShape.TextFrame.TextRange.Delete()
$titleShape.TextFrame.TextRange.Text = "A slide"

# "Just" "an" "image" "on" "this" "side" runs were merged into a single run.
$captionShape = Get-ShapeByName $s "TextBox 3"
$captionShape.TextFrame.TextRange.Delete()
$captionShape.TextFrame.TextRange.Text = "Just an image on this side"
